# Cleaned defensive actions data
# - Unmerge the grouped header cells in row 1 and give every column its own
#   header label (previously blank cells inside merged ranges).
# - Hide the old "raw" header row (row 2) and the old "totals" row (row 20);
#   insert a blank hidden spacer row (row 3).
# - Backfill missing Tkl% (column O) zeros for keepers with no tackle
#   attempts.
# - Re-point the active selection the way the workbook was left after the
#   edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: un-merge the grouped headers and label every column ---------
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2: keep the detailed header text, just hide the helper row -----
$ws.Rows.Item(2).Hidden = $true

# --- Row 3: blank hidden spacer row between the headers and the data ----
$ws.Rows.Item(3).Hidden = $true

# --- Data rows: backfill the missing Tkl% zeros --------------------------
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O19").Value = 0

# Re-save the Tkl% that had a rounding artefact so it carries full
# floating-point precision like the rest of the recomputed column.
$ws.Range("O12").Value = 33.3

# --- Row 20 (team totals) is now a hidden summary row -------------------
$ws.Rows.Item(20).Hidden = $true

# --- Selection left where the author was working when they saved --------
$ws.Range("O21").Select()
